$wb = $excel.ActiveWorkbook

# ALC!row2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 421.5
$ws.Range("J2").Value = 615
$ws.Range("L2").Value = 615
$ws.Range("N2").Value = -841

# ALC!row32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1390.3334
$ws.Range("I32").Value = 693.3333
$ws.Range("J32").Value = 1855
$ws.Range("K32").Value = 693.3333
$ws.Range("L32").Value = 1855
$ws.Range("M32").Value = -367.3333
$ws.Range("N32").Value = -2507

# ALC!row101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 315.66666
$ws.Range("J101").Value = 300
$ws.Range("L101").Value = 900
$ws.Range("N101").Value = -4144

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 3569.7
$ws.Range("I107").Value = 3099.5715
$ws.Range("J107").Value = 4666.6665
$ws.Range("K107").Value = 3099.5715
$ws.Range("L107").Value = 4666.6665
$ws.Range("M107").Value = -1179.5715
$ws.Range("N107").Value = -8506.666499999999

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3607.4
$ws.Range("I116").Value = 3044.3333
$ws.Range("J116").Value = 4452
$ws.Range("K116").Value = 3044.3333
$ws.Range("L116").Value = 4452
$ws.Range("M116").Value = 397.6667000000002
$ws.Range("N116").Value = -11336

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1624
$ws.Range("I137").Value = 1624
$ws.Range("K137").Value = 4872
$ws.Range("M137").Value = -2322

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4319.1924
$ws.Range("I32").Value = 5184.05
$ws.Range("J32").Value = 1436.3334
$ws.Range("K32").Value = 5184.05
$ws.Range("L32").Value = 1436.3334
$ws.Range("M32").Value = -4897.05
$ws.Range("N32").Value = -2010.3334

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1954.5
$ws.Range("I61").Value = 1606
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1606
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1394
$ws.Range("N61").Value = -3424

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1012.5
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 1012.5
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -516.5
$ws.Range("N97").Value = -1892

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2244.4546
$ws.Range("I110").Value = 1768
$ws.Range("J110").Value = 3078.25
$ws.Range("K110").Value = 1768
$ws.Range("L110").Value = 3078.25
$ws.Range("M110").Value = 277
$ws.Range("N110").Value = -7168.25

# ARM!row118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 38000
$ws.Range("J118").Value = 38000
$ws.Range("L118").Value = 38000
$ws.Range("N118").Value = -41314

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1609.8
$ws.Range("I122").Value = 1599.75
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 4799.25
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -2349.25
$ws.Range("N122").Value = -9850

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1954.5
$ws.Range("I136").Value = 1606
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4818
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2268
$ws.Range("N136").Value = -14100

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231670
$ws.Range("I94").Value = 25000714
$ws.Range("J94").Value = 1526.6666
$ws.Range("K94").Value = 25000714
$ws.Range("L94").Value = 1526.6666
$ws.Range("M94").Value = -25000263
$ws.Range("N94").Value = -2428.6666

# CRP!row3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4250
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 4250
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 4250
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -4476

# CRP!row20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471

# CRP!row30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181

# CRP!row128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

# CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1481.0476
$ws.Range("I5").Value = 1481.0476
$ws.Range("K5").Value = 4443.142800000001
$ws.Range("M5").Value = -4331.142800000001

# CUL!row11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1592.8
$ws.Range("I11").Value = 1866
$ws.Range("K11").Value = 5598
$ws.Range("M11").Value = -5458

# CUL!row134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4788.909
$ws.Range("I134").Value = 3451.6
$ws.Range("J134").Value = 5903.3335
$ws.Range("K134").Value = 10354.8
$ws.Range("L134").Value = 17710.0005
$ws.Range("M134").Value = -5284.799999999999
$ws.Range("N134").Value = -27850.0005

# CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1481.0476
$ws.Range("I135").Value = 1481.0476
$ws.Range("K135").Value = 13329.4284
$ws.Range("M135").Value = -10794.4284

# CUL!row136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3255.3333
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3255.3333
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 9765.999899999999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -19965.9999

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4441.25
$ws.Range("I137").Value = 1515
$ws.Range("J137").Value = 5026.5
$ws.Range("K137").Value = 4545
$ws.Range("L137").Value = 15079.5
$ws.Range("M137").Value = 555
$ws.Range("N137").Value = -25279.5

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2282.2
$ws.Range("I97").Value = 1100
$ws.Range("J97").Value = 3070.3333
$ws.Range("K97").Value = 1100
$ws.Range("L97").Value = 3070.3333
$ws.Range("M97").Value = -604
$ws.Range("N97").Value = -4062.3333

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 851.26086
$ws.Range("I16").Value = 837.1429000000001
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 837.1429000000001
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = -667.1429000000001
$ws.Range("N16").Value = -1339.5

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 896.36
$ws.Range("J22").Value = 656.0769
$ws.Range("L22").Value = 656.0769
$ws.Range("N22").Value = -1246.0769

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 896.36
$ws.Range("J27").Value = 656.0769
$ws.Range("L27").Value = 656.0769
$ws.Range("N27").Value = -870.0769

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5885.353
$ws.Range("I46").Value = 1210.2
$ws.Range("J46").Value = 7833.3335
$ws.Range("K46").Value = 1210.2
$ws.Range("L46").Value = 7833.3335
$ws.Range("M46").Value = -1022.2
$ws.Range("N46").Value = -8209.333500000001

# LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I55").Value = 98.666664
$ws.Range("J55").Value = 763.3333
$ws.Range("K55").Value = 98.666664
$ws.Range("L55").Value = 763.3333
$ws.Range("M55").Value = 74.333336
$ws.Range("N55").Value = -1109.3333

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1739.25
$ws.Range("J68").Value = 2401
$ws.Range("L68").Value = 2401
$ws.Range("N68").Value = -3899

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1739.25
$ws.Range("J71").Value = 2401
$ws.Range("L71").Value = 12005
$ws.Range("N71").Value = -19493

# WVR!row133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 38400
$ws.Range("J133").Value = 38400
$ws.Range("L133").Value = 38400
$ws.Range("N133").Value = -48520

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 841.46155
$ws.Range("I136").Value = 448.1111
$ws.Range("J136").Value = 1726.5
$ws.Range("K136").Value = 1344.3333
$ws.Range("L136").Value = 5179.5
$ws.Range("M136").Value = 1205.6667
$ws.Range("N136").Value = -10279.5
